$d = $word.ActiveDocument

# 1. "We measure and compare student" -> "We measured and compared student"
#    (bookmark _GoBack, which currently sits between "acro"/"ss " further down,
#    is untouched by this call since the match is entirely before it)
$d.Content.Find.Execute("We measure and compare student", $true, $false, $false, $false, $false,
                         $true, 1, $false, "We measured and compared student", 2)

# 2. "grades acro" + [[_GoBack bookmark]] + "ss " -> "grades across "
#    (the match text spans across the bookmark in the flattened text stream,
#    so the bookmark is consumed/removed by this replace)
$d.Content.Find.Execute("grades across", $true, $false, $false, $false, $false,
                         $true, 1, $false, "grades across", 2)

# 3. Re-create the _GoBack bookmark at its new location: right after
#    "...and compared" and right before " student efficacy...".
$text = $d.Content.Text
$pos = $text.IndexOf("We measured and compared") + ("We measured and compared").Length
$bmRange = $d.Range($pos, $pos)
$d.Bookmarks.Add("_GoBack", $bmRange)
